$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update all "Em avaliação" status cells (column C) to "Aprovada para resolução"
$ws.Range("C2:C5").Value = "Aprovada para resolução"

# Update "Data da última modificação" (column D) to 09/06/2015 (serial 42164) for rows 2-5
$ws.Range("D2").Value = [DateTime]"2015-06-09"
$ws.Range("D3").Value = [DateTime]"2015-06-09"
$ws.Range("D4").Value = [DateTime]"2015-06-09"
$ws.Range("D5").Value = [DateTime]"2015-06-09"

# Update the active selection to D5
$ws.Range("D5").Select()
